$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2210884353741497
$ws.Range("C2").Value = 0.5136054421768708
$ws.Range("J2").Value = 0.01360544217687075
$ws.Range("P2").Value = 0.1462585034013605
$ws.Range("S2").Value = 0.1054421768707483
$ws.Range("B3").Value = 0.006493506493506494
$ws.Range("C3").Value = 0.02597402597402598
$ws.Range("J3").Value = 0.03896103896103896
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2142857142857143
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.05116279069767442
$ws.Range("D6").Value = 0.02790697674418605
$ws.Range("F6").Value = 0.03255813953488372
$ws.Range("J6").Value = 0.2744186046511628
$ws.Range("O6").Value = 0.009302325581395349
$ws.Range("Q6").Value = 0.1441860465116279
$ws.Range("R6").Value = 0.06976744186046512
$ws.Range("S6").Value = 0.3906976744186046
$ws.Range("B7").Value = 0.09444444444444444
$ws.Range("D7").Value = 0.005555555555555556
$ws.Range("F7").Value = 0.03333333333333333
$ws.Range("J7").Value = 0.1388888888888889
$ws.Range("O7").Value = 0.005555555555555556
$ws.Range("Q7").Value = 0.1944444444444444
$ws.Range("R7").Value = 0.09444444444444444
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.09053497942386832
$ws.Range("D8").Value = 0.01234567901234568
$ws.Range("F8").Value = 0.051440329218107
$ws.Range("J8").Value = 0.1111111111111111
$ws.Range("O8").Value = 0.01234567901234568
$ws.Range("Q8").Value = 0.1851851851851852
$ws.Range("R8").Value = 0.09876543209876543
$ws.Range("S8").Value = 0.4382716049382716
$ws.Range("B9").Value = 0.09677419354838709
$ws.Range("D9").Value = 0.01209677419354839
$ws.Range("F9").Value = 0.06451612903225806
$ws.Range("J9").Value = 0.1048387096774194
$ws.Range("O9").Value = 0.008064516129032258
$ws.Range("Q9").Value = 0.1774193548387097
$ws.Range("R9").Value = 0.1048387096774194
$ws.Range("S9").Value = 0.4314516129032258
$ws.Range("B10").Value = 0.09748667174409749
$ws.Range("D10").Value = 0.02284843869002285
$ws.Range("F10").Value = 0.06702208682406702
$ws.Range("J10").Value = 0.1287128712871287
$ws.Range("O10").Value = 0.009900990099009901
$ws.Range("Q10").Value = 0.2079207920792079
$ws.Range("R10").Value = 0.09367859862909368
$ws.Range("S10").Value = 0.3724295506473724
$ws.Range("G11").Value = 0.1360294117647059
$ws.Range("J11").Value = 0.07720588235294118
$ws.Range("K11").Value = 0.1727941176470588
$ws.Range("L11").Value = 0.5955882352941176
$ws.Range("S11").Value = 0.01838235294117647
$ws.Range("G12").Value = 0.8098159509202454
$ws.Range("J12").Value = 0.1656441717791411
$ws.Range("K12").Value = 0.006134969325153374
$ws.Range("L12").Value = 0.01226993865030675
$ws.Range("S12").Value = 0.006134969325153374
$ws.Range("G13").Value = 0.5294117647058824
$ws.Range("J13").Value = 0.3823529411764706
$ws.Range("S13").Value = 0.08823529411764706
$ws.Range("F15").Value = 0.02538071065989848
$ws.Range("H15").Value = 0.1776649746192893
$ws.Range("I15").Value = 0.1015228426395939
$ws.Range("J15").Value = 0.3096446700507614
$ws.Range("K15").Value = 0.1015228426395939
$ws.Range("M15").Value = 0.02030456852791878
$ws.Range("O15").Value = 0.04568527918781726
$ws.Range("S15").Value = 0.2182741116751269
$ws.Range("F16").Value = 0.01621621621621622
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.1027027027027027
$ws.Range("J16").Value = 0.3837837837837838
$ws.Range("K16").Value = 0.0918918918918919
$ws.Range("M16").Value = 0.01081081081081081
$ws.Range("O16").Value = 0.06486486486486487
$ws.Range("S16").Value = 0.1297297297297297
$ws.Range("F17").Value = 0.01923076923076923
$ws.Range("H17").Value = 0.1602564102564103
$ws.Range("I17").Value = 0.1025641025641026
$ws.Range("J17").Value = 0.4551282051282051
$ws.Range("K17").Value = 0.06623931623931624
$ws.Range("M17").Value = 0.01495726495726496
$ws.Range("N17").Value = 0.002136752136752137
$ws.Range("O17").Value = 0.06837606837606838
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("F18").Value = 0.01304347826086956
$ws.Range("H18").Value = 0.1956521739130435
$ws.Range("I18").Value = 0.09565217391304348
$ws.Range("J18").Value = 0.4478260869565218
$ws.Range("K18").Value = 0.1043478260869565
$ws.Range("M18").Value = 0.01739130434782609
$ws.Range("O18").Value = 0.02173913043478261
$ws.Range("S18").Value = 0.1043478260869565
$ws.Range("F19").Value = 0.01650412603150788
$ws.Range("H19").Value = 0.2273068267066767
$ws.Range("I19").Value = 0.1042760690172543
$ws.Range("J19").Value = 0.3563390847711928
$ws.Range("K19").Value = 0.09452363090772693
$ws.Range("M19").Value = 0.01575393848462115
$ws.Range("N19").Value = 0.0007501875468867217
$ws.Range("O19").Value = 0.06601650412603151
$ws.Range("S19").Value = 0.118529632408102
